$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 58. Excel shifts rows 58-60 down to 59-61 and
# copies row 58's formatting into the freshly inserted blank row.
$ws.Rows(58).Insert()

# Populate the newly inserted row 58 with the new lot number (27720) and
# the same standard-curve-id columns used by the surrounding rows.
$ws.Range("A58").Value = 27720
$ws.Range("B58").Value = "GX23-0002_SARS_COV_2"
$ws.Range("C58").Value = "GX23-0001_FluA1"
$ws.Range("D58").Value = "GX23-0001_FluA2"
$ws.Range("E58").Value = "GX23-0001_FluB"
$ws.Range("F58").Value = "GX23-0001_RSV"

# Append two new rows (62, 63) for the additional lot numbers 48116 and
# 49801, re-using the same standard-curve-id values as the rows above.
$ws.Range("A62").Value = 48116
$ws.Range("B62").Value = "GX23-0002_SARS_COV_2"
$ws.Range("C62").Value = "GX23-0001_FluA1"
$ws.Range("D62").Value = "GX23-0001_FluA2"
$ws.Range("E62").Value = "GX23-0001_FluB"
$ws.Range("F62").Value = "GX23-0001_RSV"

$ws.Range("A63").Value = 49801
$ws.Range("B63").Value = "GX23-0002_SARS_COV_2"
$ws.Range("C63").Value = "GX23-0001_FluA1"
$ws.Range("D63").Value = "GX23-0001_FluA2"
$ws.Range("E63").Value = "GX23-0001_FluB"
$ws.Range("F63").Value = "GX23-0001_RSV"

# All of the lot-number cells (A58:A63) pick up the new "most recent std
# curve" formatting: left-aligned, top-vertical-aligned.
$lotRange = $ws.Range("A58:A63")
$lotRange.HorizontalAlignment = -4131
$lotRange.VerticalAlignment = -4160

# Move the active selection, matching the saved workbook state.
$ws.Range("C55").Select() | Out-Null
